{"js": "// Grading section update:\n//   1. \"Assignments (60%)\" -> \"Assignments (50%)\"\n//   2. New list item \"Exams / Quizzes (10%)\" added after \"Labs (15%)\"\n//      (and before \"Final (10%)\"), matching the existing bullet/numbering\n//      list used by the other grading line items.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nlet assignmentsPara = null;\nlet labsPara = null;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  const t = paragraphs.items[i].text;\n  if (t === \"Assignments (60%)\") {\n    assignmentsPara = paragraphs.items[i];\n  } else if (t === \"Labs (15%)\") {\n    labsPara = paragraphs.items[i];\n  }\n}\n\nif (!assignmentsPara) {\n  throw new Error('Could not find the \"Assignments (60%)\" paragraph.');\n}\nif (!labsPara) {\n  throw new Error('Could not find the \"Labs (15%)\" paragraph.');\n}\n\n// 1) Update the Assignments weight from 60% to 50%.\nconst pctMatches = assignmentsPara.getRange().search(\"60\", { matchCase: true });\npctMatches.load(\"items\");\nawait context.sync();\nif (pctMatches.items.length === 0) {\n  throw new Error('Could not find \"60\" inside the Assignments paragraph.');\n}\npctMatches.items[0].insertText(\"50\", \"Replace\");\n\n// 2) Insert a new grading line item for quizzes, right after \"Labs (15%)\",\n//    re-using the same list (numbering) as the surrounding bullet items.\nconst list = labsPara.list;\nlist.load(\"id\");\nawait context.sync();\n\nconst quizzesPara = labsPara.insertParagraph(\"Exams / Quizzes (10%)\", \"After\");\nquizzesPara.style = \"List Paragraph\";\nquizzesPara.attachToList(list.id, 0);\n\nawait context.sync();\n", "ps1": "# Grading section update:\n#   1. \"Assignments (60%)\" -> \"Assignments (50%)\"\n#   2. New list item \"Exams / Quizzes (10%)\" added after \"Labs (15%)\"\n#      (and before \"Final (10%)\"), matching the existing bullet/numbering\n#      list used by the other grading line items.\n\n$d = $word.ActiveDocument\n$paras = $d.Paragraphs\n$count = $paras.Count\n\n$assignmentsPara = $null\n$labsPara = $null\nfor ($i = 1; $i -le $count; $i++) {\n    $p = $paras.Item($i)\n    $t = $p.Range.Text\n    if ($t -eq \"Assignments (60%)`r\") {\n        $assignmentsPara = $p\n    } elseif ($t -eq \"Labs (15%)`r\") {\n        $labsPara = $p\n    }\n}\n\nif ($null -eq $assignmentsPara) {\n    throw \"Could not find the 'Assignments (60%)' paragraph.\"\n}\nif ($null -eq $labsPara) {\n    throw \"Could not find the 'Labs (15%)' paragraph.\"\n}\n\n# 1) Update the Assignments weight from 60% to 50%, scoped to just that\n#    paragraph so the \"60.00 - 68.99\" grade-scale table row is untouched.\n$find = $assignmentsPara.Range.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"60\"\n$find.Replacement.Text = \"50\"\n$find.Execute($find.Text, $false, $true, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 1) | Out-Null\n\n# 2) Insert a new grading line item for quizzes, right after \"Labs (15%)\".\n#    InsertParagraphAfter copies the paragraph formatting (ListParagraph\n#    style + numbering) from \"Labs (15%)\" onto the new paragraph.\n$labsPara.Range.InsertParagraphAfter()\n\n$paras2 = $d.Paragraphs\n$count2 = $paras2.Count\nfor ($i = 1; $i -le $count2; $i++) {\n    $p = $paras2.Item($i)\n    if ($p.Range.Text -eq \"Labs (15%)`r\") {\n        $newPara = $paras2.Item($i + 1)\n        $newPara.Range.Text = \"Exams / Quizzes (10%)\"\n        break\n    }\n}\n"}
